$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999088063629
$ws.Range("A2").Value = 0.99770837692795289
$ws.Range("A3").Value = 0.9935591722672299
$ws.Range("A4").Value = 0.99638236753418241
$ws.Range("A5").Value = 0.98605844551271749
$ws.Range("A6").Value = 0.96089296566472138
$ws.Range("A7").Value = 0.95902581401898579
$ws.Range("A8").Value = 0.95636660130683993
$ws.Range("A9").Value = 0.95607070031746821
$ws.Range("A10").Value = 0.95672585775126739
$ws.Range("A11").Value = 0.95689095427797266
$ws.Range("A12").Value = 0.95396426572436732
$ws.Range("A13").Value = 0.94146110044102915
$ws.Range("A14").Value = 0.93729441910923361
$ws.Range("A15").Value = 0.93470327566550782
$ws.Range("A16").Value = 0.93219701021300549
$ws.Range("A17").Value = 0.92848932290977215
$ws.Range("A18").Value = 0.9273804345545672
$ws.Range("A19").Value = 0.99467086806764016
$ws.Range("A20").Value = 0.97004106670252155
$ws.Range("A21").Value = 0.96321178949893138
$ws.Range("A22").Value = 0.96194727767385602
$ws.Range("A23").Value = 0.98652505560170378
$ws.Range("A24").Value = 0.97350485942562637
$ws.Range("A25").Value = 0.96704798369727407
$ws.Range("A26").Value = 0.95723967799557186
$ws.Range("A27").Value = 0.95239585804308191
$ws.Range("A28").Value = 0.93093198131008736
$ws.Range("A29").Value = 0.91566402778262623
$ws.Range("A30").Value = 0.90909438152065103
$ws.Range("A31").Value = 0.90144127992589074
$ws.Range("A32").Value = 0.89976202031928709
$ws.Range("A33").Value = 0.89924202562772071
